$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)

# Resize the text box
$shp.Width = 8873711 / 914400 * 72
$shp.Height = 3893374 / 914400 * 72

$tr = $shp.TextFrame.TextRange

# Insert a new first paragraph "Exploratory Data Analysis" before the
# existing "G2M Case Study" paragraph, inheriting its formatting.
$para1 = $tr.Paragraphs(1,1)
[void]$para1.InsertBefore("Exploratory Data Analysis`r")

# Update the date text (edit the run directly so the single run is preserved)
$paras = $tr.Paragraphs()
foreach ($para in $paras) {
    if ($para.Text -eq "20-Jan-2021") {
        $run = $para.Runs(1,1)
        $run.Text = "20-Jan-2023"
    }
}
